# Generate Report for Handoff
# -----------------------------------------------------------------------
# A second source file (ffff19a36d12-436c-4a83-a9a2-4c2d3b2a046f.md) has
# been handed off together with the existing file - whose own uuid/hash
# also rolled to a new handoff (e7755907-0651-44dc-b0ea-594ad263c848 /
# 1c3e61117fa6e114a9ce0cd50106b4b01b066a13). This inserts a row for the
# new file (sharing the same handoff batch/timestamps as the first file)
# ahead of the existing ".localization-config" row on every sheet.

$wb = $excel.ActiveWorkbook

$oldMd      = "135b7efa-90b3-44f6-b424-ee0e50e8548d.md"
$newMd1     = "e7755907-0651-44dc-b0ea-594ad263c848.md"
$newMd2     = "ffff19a36d12-436c-4a83-a9a2-4c2d3b2a046f.md"
$config     = ".localization-config"

$xlfZh      = "e7755907-0651-44dc-b0ea-594ad263c848.1c3e61117fa6e114a9ce0cd50106b4b01b066a13.zh-cn.xlf"
$xlfDe      = "e7755907-0651-44dc-b0ea-594ad263c848.1c3e61117fa6e114a9ce0cd50106b4b01b066a13.de-de.xlf"

$dtZh       = "2016-03-03 10:47:46"
$dtDe       = "2016-03-03 10:47:58"
$dtEmpty    = "0001-01-01 00:00:00"

$urlRepoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/a114a5983bd14b6750400c67666f0e1acf2e2db5"
$urlMd1        = "$urlRepoBase/e2e/$newMd1"
$urlMd2        = "$urlRepoBase/e2e/$newMd2"
$urlConfig     = "$urlRepoBase/$config"
$urlXlfZh      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60679bb21db82b6ca50183a01e493da5bf1e0825/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh"
$urlXlfDe      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d6b78102d25a8e1fe572ba84069d5a6b01b017e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe"

# ---------------------------------------------------------------------
# Sheet "Overview" - 3 columns (File Name / zh-cn / de-de)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd1, "", "", $newMd1)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, "", "", $newMd2)
$ws.Hyperlinks.Add($ws.Range("A4"), $urlConfig, "", "", $config)

# ---------------------------------------------------------------------
# Sheet "zh-cn" - 9 columns (Source File Name ... Dependency From)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $xlfZh
$ws.Range("D2").Value = $dtZh
$ws.Range("G2").Value = $dtEmpty
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $xlfZh
$ws.Range("D3").Value = $dtZh
$ws.Range("G3").Value = $dtEmpty
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = $dtEmpty
$ws.Range("G4").Value = $dtEmpty
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd1, "", "", $newMd1)
$ws.Hyperlinks.Add($ws.Range("C2"), $urlXlfZh, "", "", $xlfZh)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, "", "", $newMd2)
$ws.Hyperlinks.Add($ws.Range("C3"), $urlXlfZh, "", "", $xlfZh)
$ws.Hyperlinks.Add($ws.Range("A4"), $urlConfig, "", "", $config)

# ---------------------------------------------------------------------
# Sheet "de-de" - 9 columns (Source File Name ... Dependency From)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $xlfDe
$ws.Range("D2").Value = $dtDe
$ws.Range("G2").Value = $dtEmpty
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $xlfDe
$ws.Range("D3").Value = $dtDe
$ws.Range("G3").Value = $dtEmpty
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = $dtEmpty
$ws.Range("G4").Value = $dtEmpty
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd1, "", "", $newMd1)
$ws.Hyperlinks.Add($ws.Range("C2"), $urlXlfDe, "", "", $xlfDe)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, "", "", $newMd2)
$ws.Hyperlinks.Add($ws.Range("C3"), $urlXlfDe, "", "", $xlfDe)
$ws.Hyperlinks.Add($ws.Range("A4"), $urlConfig, "", "", $config)
